$d = $word.ActiveDocument

# Collapse to the very end of the document body (after the last existing
# paragraph, before the section break) so the new paragraphs are appended
# after "This device will use a microcontroller ..." paragraph.
$r = $d.Content
$r.Collapse(0)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>It will be a vertically</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">oriented rectangle </w:t></w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:t>similar to</w:t></w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:t xml:space="preserve"> an identification badge. The main rectangle will be 3 in. wide and 4 in. tall with rounded corners.</w:t></w:r>
          </w:p>
          <w:p>
            <w:r><w:t>It will have an attachment point to be connected to a lanyard.</w:t></w:r>
          </w:p>
          <w:p>
            <w:r><w:t>It will be powered by a USB-rechargeable battery using a USB-C connector</w:t></w:r>
            <w:r><w:t xml:space="preserve"> (if possible, dual-purpose the connector for reprogramming the microcontroller)</w:t></w:r>
            <w:r><w:t xml:space="preserve">. It will have a piezoelectric buzzer for playing sounds. </w:t></w:r>
            <w:r><w:t>It can also connect to a heart rate monitor, if possible.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
